$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column headers I1 ("I0") and J1 ("IF"), matching the style of
# the existing header cells (e.g. H1) which use bold/centered/bordered style.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Populate the new I (I0) and J (IF) columns for data rows 2-69.
$data = @(
    @(2,9,9),
    @(3,8,8),
    @(4,5,5),
    @(5,5,5),
    @(6,5,5),
    @(7,8,8),
    @(8,5,5),
    @(9,9,9),
    @(10,8,8),
    @(11,7,7),
    @(12,9,9),
    @(13,8,8),
    @(14,8,8),
    @(15,6,6),
    @(16,8,8),
    @(17,9,9),
    @(18,9,9),
    @(19,8,8),
    @(20,8,8),
    @(21,8,8),
    @(22,7,7),
    @(23,9,9),
    @(24,7,7),
    @(25,8,8),
    @(26,8,8),
    @(27,8,8),
    @(28,7,7),
    @(29,9,9),
    @(30,8,8),
    @(31,8,8),
    @(32,9,9),
    @(33,9,9),
    @(34,8,8),
    @(35,9,9),
    @(36,9,9),
    @(37,9,9),
    @(38,8,8),
    @(39,8,8),
    @(40,7,7),
    @(41,7,7),
    @(42,8,8),
    @(43,8,8),
    @(44,7,7),
    @(45,8,8),
    @(46,7,7),
    @(47,6,6),
    @(48,7,7),
    @(49,7,7),
    @(50,6,7),
    @(51,6,6),
    @(52,6,6),
    @(53,8,8),
    @(54,10,10),
    @(55,6,6),
    @(56,7,7),
    @(57,8,8),
    @(58,9,9),
    @(59,6,6),
    @(60,7,7),
    @(61,7,7),
    @(62,9,9),
    @(63,7,7),
    @(64,6,6),
    @(65,3,3),
    @(66,5,5),
    @(67,6,6),
    @(68,6,6),
    @(69,3,3)
)

foreach ($row in $data) {
    $r = $row[0]
    $i = $row[1]
    $j = $row[2]
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $j
}
